$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the existing "Lato" cell formatting (currently on B2) and
#        stamp it onto the two new header cells (A2, A10) before we touch
#        B2's own formatting. PasteSpecial(formats) only copies the
#        format, not the value, so doing this first is safe.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Reset B2/A3/A7 back to the plain/default look: B2 no longer
#        carries the Lato styling, and A3/A7 no longer carry the
#        Hyperlink styling in the new layout. Copy A1's plain format over.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Rewrite row 2 / row 3 (CarLoanAmount block) -----------------------
$ws.Range("A2").Value = "CarLoanAmount"
$ws.Range("B2").Value = "InterestRate"
$ws.Range("C2").Value = "LoanTenure"
$ws.Range("D2").ClearContents()

$ws.Range("A3").Value = 1500000
$ws.Range("B3").Value = 9.5
$ws.Range("C3").Value = 1
$ws.Range("D3").ClearContents()

# --- 4. Rewrite the title block (rows 5-7) --------------------------------
$ws.Range("A5").Value = "driverTitleTestData"

$ws.Range("A6").Value = "title"
$ws.Range("B6").ClearContents()

$ws.Range("A7").Value = "EMI Calculator for Home Loan, Car Loan & Personal Loan in India"
$ws.Range("B7").ClearContents()

# --- 5. New smoke-suite block: totalPaymentAmountVerify (rows 9-11) ------
$ws.Range("A9").Value = "totalPaymentAmountVerify"

$ws.Range("A10").Value = "CarLoanAmount"
$ws.Range("B10").Value = "InterestRate"
$ws.Range("C10").Value = "LoanTenure"
$ws.Range("D10").Value = "expectedTotalAmount"

$ws.Range("A11").Value = 1000000
$ws.Range("B11").Value = 12.5
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 1349876
$ws.Range("D11").NumberFormat = "#,##0"

# --- 6. Drop the hyperlinks + the now-unused Hyperlink cell style --------
$ws.Hyperlinks.Delete()
$wb.Styles.Item("Hyperlink").Delete()

# --- 7. Column D needs to be a bit wider for "expectedTotalAmount" -------
$ws.Columns("D").ColumnWidth = 17

# --- 8. Final selection / active cell -------------------------------------
$ws.Range("A11").Select()
